# feat: add 2022-Q3 data
#
# The workbook currently has 4 sheets: 总计, 2022-Q2, 2022-Q1, 2021-Q4.
# We insert a brand-new "2022-Q3" sheet right after "总计" (so it becomes
# the 2nd tab, pushing 2022-Q2/2022-Q1/2021-Q4 one slot to the right),
# populate it with the Q3 fund-holding table, and insert a matching
# summary row at the top of the "总计" sheet's data.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q3" worksheet right after "总计".
# ---------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q3Sheet.Name = "2022-Q3"

# NOTE: fetch this handle *after* the Add()/rename above — a handle grabbed
# beforehand tracks the worksheet *position*, and Add() shifts "2022-Q2"
# over by one slot, so a pre-fetched reference would silently resolve to
# the freshly inserted "2022-Q3" sheet instead.
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# Reuse the header / index-column formatting from the 2022-Q2 sheet so the
# new sheet's look (bold+bordered header row, bold+bordered column A) is
# consistent with the rest of the workbook. The 2022-Q2 table has exactly
# the 4 rows (1 header + 3 data rows) we need for Q3 as well. Column A1
# itself is always blank in these tables, so copy the header (B1:H1) and
# the index column (A2:A4) separately rather than the full A1:H4 block.
$q2Sheet.Range("B1:H1").Copy()
$q3Sheet.Range("B1").PasteSpecial(-4122)   # xlPasteFormats
$q2Sheet.Range("A2:A4").Copy()
$q3Sheet.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

$q3Headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $q3Headers.Length; $i++) {
    $q3Sheet.Cells.Item(1, 2 + $i).Value = $q3Headers[$i]
}

# index, 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$q3Data = @(
    @(0, "851088", "海通量化成长精选一年持有期混合A", "0.38", "85.56", "0.92", "0.0035", 7),
    @(1, "850010", "海通量化成长精选一年持有期混合B", "0.25", "85.56", "0.92", "0.0023", 7),
    @(2, "851099", "海通量化成长精选一年持有期混合C", "0.03", "85.56", "0.92", "0.0003", 7)
)

for ($r = 0; $r -lt $q3Data.Length; $r++) {
    $row = $q3Data[$r]
    $excelRow = 2 + $r
    $q3Sheet.Cells.Item($excelRow, 1).Value = $row[0]
    # Leading apostrophe forces text storage (fund codes / ratios are text
    # in the source data, not numbers) while keeping the actual cell value
    # free of the apostrophe itself. Re-flatten the style afterwards back
    # to Normal since the quote-prefix write bumps the cell onto a new
    # style record otherwise.
    $q3Sheet.Cells.Item($excelRow, 2).Value = "'" + $row[1]
    $q3Sheet.Cells.Item($excelRow, 2).Style = "Normal"
    $q3Sheet.Cells.Item($excelRow, 3).Value = $row[2]
    $q3Sheet.Cells.Item($excelRow, 4).Value = "'" + $row[3]
    $q3Sheet.Cells.Item($excelRow, 4).Style = "Normal"
    $q3Sheet.Cells.Item($excelRow, 5).Value = "'" + $row[4]
    $q3Sheet.Cells.Item($excelRow, 5).Style = "Normal"
    $q3Sheet.Cells.Item($excelRow, 6).Value = "'" + $row[5]
    $q3Sheet.Cells.Item($excelRow, 6).Style = "Normal"
    $q3Sheet.Cells.Item($excelRow, 7).Value = "'" + $row[6]
    $q3Sheet.Cells.Item($excelRow, 7).Style = "Normal"
    $q3Sheet.Cells.Item($excelRow, 8).Value = $row[7]
}

# ---------------------------------------------------------------------
# 2) Insert the new Q3 summary row at the top of the "总计" table's data
#    (row 2), pushing the existing 2022-Q2 / 2022-Q1 / 2021-Q4 rows down.
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

# Drop the formatting the Insert() copied down from the header row above
# for the plain data cells (B:D), then restore column A's index style by
# copying it from the row just below (which still has the original style).
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 3
$totalSheet.Cells.Item(2, 4).Value = 0.01

# The A column is a simple 0-based row counter, independent of the quarter
# label; renumber the rows that got pushed down by the insert.
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3
